$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column letter -> new value (G:T), only columns that change.
$updates = @{
    2  = @{ G='26.81310566666667'; H='80.439317'; I='0.004518206005002021'; J='0.004518206005002021'; K='3'; L='1'; M='1.728078666666667'; N='5.184236'; O='0.01495566191401176'; P='0.01495566191401175'; Q='46.33515588964578'; R='417.016403006812'; S='6.757276146866793E-05'; T='6.757276146866792E-05' }
    3  = @{ G='26.81310566666667'; H='80.439317'; I='0.004518206005002021'; J='0.004518206005002021'; O='0.6765388985664503'; P='0.6765388985664502'; Q='2096.031289736274'; R='18864.28160762646'; S='0.003056742114120389'; T='0.003056742114120388' }
    4  = @{ G='26.81310566666667'; H='80.439317'; I='0.004518206005002021'; J='0.004518206005002021'; M='35.284818'; N='105.854454'; O='0.3053725613795956'; P='0.3053725613795956'; Q='946.0955534631021'; R='8514.859981167918'; S='0.001379736140588137'; T='0.001379736140588137' }
    5  = @{ G='26.81310566666667'; H='80.439317'; I='0.004518206005002021'; J='0.004518206005002021'; K='3'; L='1'; M='0.361994'; N='1.085982'; O='0.003132878139942378'; P='0.003132878139942378'; Q='9.706183372699334'; R='87.355650354294'; S='1.415498882482721E-05'; T='1.415498882482721E-05' }
    6  = @{ G='5771.873535333333'; I='0.9726032482643521'; J='0.9726032482643523'; K='3'; L='1'; M='1.728078666666667'; N='5.184236'; O='0.01495566191401176'; P='0.01495566191401175'; Q='9974.251523107447'; R='89768.26370796702'; S='0.01454592535751129'; T='0.01454592535751129' }
    7  = @{ G='5771.873535333333'; I='0.9726032482643521'; J='0.9726032482643523'; O='0.6765388985664503'; P='0.6765388985664502'; Q='451198.293881314'; R='4060784.644931826'; S='0.6580039303229166'; T='0.6580039303229166' }
    8  = @{ G='5771.873535333333'; I='0.9726032482643521'; J='0.9726032482643523'; M='35.284818'; N='105.854454'; O='0.3053725613795956'; P='0.3053725613795956'; Q='203659.5072132532'; R='1832935.564919279'; S='0.2970063451285999'; T='0.2970063451286' }
    9  = @{ G='5771.873535333333'; I='0.9726032482643521'; J='0.9726032482643523'; K='3'; L='1'; M='0.361994'; N='1.085982'; O='0.003132878139942378'; P='0.003132878139942378'; Q='2089.383588549455'; R='18804.45229694509'; S='0.003047047455324338'; T='0.003047047455324338' }
    10 = @{ G='132.4457753333333'; H='397.337326'; I='0.02231808970163987'; J='0.02231808970163988'; K='3'; L='1'; M='1.728078666666667'; N='5.184236'; O='0.01495566191401176'; P='0.01495566191401175'; Q='228.8767188436596'; R='2059.890469592936'; S='0.0003337818041443134'; T='0.0003337818041443135' }
    11 = @{ G='132.4457753333333'; H='397.337326'; I='0.02231808970163987'; J='0.02231808970163988'; O='0.6765388985664503'; P='0.6765388985664502'; Q='10353.53728669951'; R='93181.83558029565'; S='0.01509905582485468'; T='0.01509905582485468' }
    12 = @{ G='132.4457753333333'; H='397.337326'; I='0.02231808970163987'; J='0.02231808970163988'; M='35.284818'; N='105.854454'; O='0.3053725613795956'; P='0.3053725613795956'; Q='4673.325077505556'; R='42059.92569755001'; S='0.006815332217289343'; T='0.006815332217289343' }
    13 = @{ G='132.4457753333333'; H='397.337326'; I='0.02231808970163987'; J='0.02231808970163988'; K='3'; L='1'; M='0.361994'; N='1.085982'; O='0.003132878139942378'; P='0.003132878139942378'; Q='47.94457599601466'; R='431.501183964132'; S='6.991985535154067E-05'; T='6.991985535154067E-05' }
    14 = @{ G='3.326003'; H='9.978009'; I='0.0005604560290058679'; J='0.000560456029005868'; K='3'; L='1'; M='1.728078666666667'; N='5.184236'; O='0.01495566191401176'; P='0.01495566191401175'; Q='5.747594829569334'; R='51.72835346612401'; S='8.381990887481327E-06'; T='8.381990887481327E-06' }
    15 = @{ G='3.326003'; H='9.978009'; I='0.0005604560290058679'; J='0.000560456029005868'; O='0.6765388985664503'; P='0.6765388985664502'; Q='259.9999583943527'; R='2339.999625549174'; S='0.0003791703045585564'; T='0.0003791703045585564' }
    16 = @{ G='3.326003'; H='9.978009'; I='0.0005604560290058679'; J='0.000560456029005868'; M='35.284818'; N='105.854454'; O='0.3053725613795956'; P='0.3053725613795956'; Q='117.357410522454'; R='1056.216694702086'; S='0.0001711478931181588'; T='0.0001711478931181588' }
    17 = @{ G='3.326003'; H='9.978009'; I='0.0005604560290058679'; J='0.000560456029005868'; K='3'; L='1'; M='0.361994'; N='1.085982'; O='0.003132878139942378'; P='0.003132878139942378'; Q='1.203993129982'; R='10.835938169838'; S='1.755840441671395E-06'; T='1.755840441671395E-06' }
}

foreach ($rowNum in $updates.Keys) {
    $rowVals = $updates[$rowNum]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$rowNum").Value = [double]$rowVals[$col]
    }
}
